# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# Rewrites the worker/Periodo Mora detail table (rows 16-34) on Hoja1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, TipoDoc(B), NumDoc(C), Nombre(D), PeriodoMora(E), ValorMora(F), SalarioBasico(G)
$rows = @(
    @{R=16; B="CC"; C="22800993";   D="DEIVIS DE JESUS PAJARO ESPINOSA";    E="1903"; F=33054; G=826361},
    @{R=17; B="CC"; C="45514619";   D="LIDIA MARIA CARDONA PEREZ";          E="1903"; F=33054; G=826361},
    @{R=18; B="CC"; C="1050945650"; D="ANGELINA ROSSI GONZALEZ";            E="1903"; F=66000; G=1650000},
    @{R=19; B="CC"; C="1050954580"; D="JEINNY ALEJANDRA DE LA CRUZ PUELLO"; E="1612"; F=28092; G=702000},
    @{R=20; B="CC"; C="1128049796"; D="CLAUDIA RAFAELA PADILLA RAMIREZ";    E="1911"; F=28647; G=828116},
    @{R=21; B="CC"; C="1128049796"; D="CLAUDIA RAFAELA PADILLA RAMIREZ";    E="1910"; F=33054; G=828116},
    @{R=22; B="CC"; C="1128049796"; D="CLAUDIA RAFAELA PADILLA RAMIREZ";    E="1909"; F=33054; G=828116},
    @{R=23; B="CC"; C="1128049796"; D="CLAUDIA RAFAELA PADILLA RAMIREZ";    E="1908"; F=33054; G=828116},
    @{R=24; B="CC"; C="1128049796"; D="CLAUDIA RAFAELA PADILLA RAMIREZ";    E="1907"; F=33054; G=828116},
    @{R=25; B="CC"; C="1128049796"; D="CLAUDIA RAFAELA PADILLA RAMIREZ";    E="1906"; F=33054; G=828116},
    @{R=26; B="CC"; C="64702479";   D="GISELA LORENA LOPEZ ENAMORADO";      E="2406"; F=17333; G=1300000},
    @{R=27; B="CC"; C="1044906744"; D="MARIA CLAUDIA JAIMES PEREIRA";       E="1911"; F=1104;  G=828116},
    @{R=28; B="CC"; C="1044906744"; D="MARIA CLAUDIA JAIMES PEREIRA";       E="1910"; F=33125; G=828116},
    @{R=29; B="CC"; C="1044906744"; D="MARIA CLAUDIA JAIMES PEREIRA";       E="1909"; F=33125; G=828116},
    @{R=30; B="CC"; C="1044906744"; D="MARIA CLAUDIA JAIMES PEREIRA";       E="1908"; F=17667; G=828116},
    @{R=31; B="CC"; C="20325192";   D="MARIA MARGARITA GARAY JIMENEZ";      E="2109"; F=33600; G=1800000},
    @{R=32; B="CC"; C="1193210960"; D="SEBASTIAN JOSE BOLAÑO COGOLLO";      E="2110"; F=27066; G=1000000},
    @{R=33; B="CC"; C="1143346176"; D="EVIS ADRIANA GUERRA BOLIVAR";        E="2403"; F=72000; G=1800000},
    @{R=34; B="CC"; C="30773549";   D="ROCIO DEL CARMEN PASSO ALVIS";       E="2406"; F=27733; G=1300000}
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
}
